$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.558564
$ws.Range("H2").Value = 4.675692
$ws.Range("I2").Value = 0.005692101168584756
$ws.Range("J2").Value = 0.005692101168584756
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.397241999999999
$ws.Range("N2").Value = 16.191726
$ws.Range("O2").Value = 0.4331003391330405
$ws.Range("P2").Value = 0.4331003391330406
$ws.Range("Q2").Value = 8.411947080487998
$ws.Range("R2").Value = 75.70752372439199
$ws.Range("S2").Value = 0.002465250946493634
$ws.Range("T2").Value = 0.002465250946493634

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.558564
$ws.Range("H3").Value = 4.675692
$ws.Range("I3").Value = 0.005692101168584756
$ws.Range("J3").Value = 0.005692101168584756
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.372979
$ws.Range("N3").Value = 13.118937
$ws.Range("O3").Value = 0.3509086099755513
$ws.Range("P3").Value = 0.3509086099755513
$ws.Range("Q3").Value = 6.815567642155999
$ws.Range("R3").Value = 61.34010877940399
$ws.Range("S3").Value = 0.001997407308908288
$ws.Range("T3").Value = 0.001997407308908288

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.558564
$ws.Range("H4").Value = 4.675692
$ws.Range("I4").Value = 0.005692101168584756
$ws.Range("J4").Value = 0.005692101168584756
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.59901
$ws.Range("N4").Value = 1.79703
$ws.Range("O4").Value = 0.04806740815847847
$ws.Range("P4").Value = 0.04806740815847847
$ws.Range("Q4").Value = 0.93359542164
$ws.Range("R4").Value = 8.40235879476
$ws.Range("S4").Value = 0.0002736045501497157
$ws.Range("T4").Value = 0.0002736045501497158

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.558564
$ws.Range("H5").Value = 4.675692
$ws.Range("I5").Value = 0.005692101168584756
$ws.Range("J5").Value = 0.005692101168584756
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.092643333333333
$ws.Range("N5").Value = 6.27793
$ws.Range("O5").Value = 0.1679236427329297
$ws.Range("P5").Value = 0.1679236427329297
$ws.Range("Q5").Value = 3.261518564173333
$ws.Range("R5").Value = 29.35366707756
$ws.Range("S5").Value = 0.0009558383630331184
$ws.Range("T5").Value = 0.0009558383630331184

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 259.5505726666667
$ws.Range("H6").Value = 778.6517180000001
$ws.Range("I6").Value = 0.9479162344201305
$ws.Range("J6").Value = 0.9479162344201304
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.397241999999999
$ws.Range("N6").Value = 16.191726
$ws.Range("O6").Value = 0.4331003391330405
$ws.Range("P6").Value = 0.4331003391330406
$ws.Range("Q6").Value = 1400.857251920585
$ws.Range("R6").Value = 12607.71526728527
$ws.Range("S6").Value = 0.4105428425970732
$ws.Range("T6").Value = 0.4105428425970732

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 259.5505726666667
$ws.Range("H7").Value = 778.6517180000001
$ws.Range("I7").Value = 0.9479162344201305
$ws.Range("J7").Value = 0.9479162344201304
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.372979
$ws.Range("N7").Value = 13.118937
$ws.Range("O7").Value = 0.3509086099755513
$ws.Range("P7").Value = 0.3509086099755513
$ws.Range("Q7").Value = 1135.009203709307
$ws.Range("R7").Value = 10215.08283338377
$ws.Range("S7").Value = 0.3326319681936269
$ws.Range("T7").Value = 0.3326319681936268

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 259.5505726666667
$ws.Range("H8").Value = 778.6517180000001
$ws.Range("I8").Value = 0.9479162344201305
$ws.Range("J8").Value = 0.9479162344201304
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.59901
$ws.Range("N8").Value = 1.79703
$ws.Range("O8").Value = 0.04806740815847847
$ws.Range("P8").Value = 0.04806740815847847
$ws.Range("Q8").Value = 155.47338853306
$ws.Range("R8").Value = 1399.26049679754
$ws.Range("S8").Value = 0.04556387653992037
$ws.Range("T8").Value = 0.04556387653992037

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 259.5505726666667
$ws.Range("H9").Value = 778.6517180000001
$ws.Range("I9").Value = 0.9479162344201305
$ws.Range("J9").Value = 0.9479162344201304
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.092643333333333
$ws.Range("N9").Value = 6.27793
$ws.Range("O9").Value = 0.1679236427329297
$ws.Range("P9").Value = 0.1679236427329297
$ws.Range("Q9").Value = 543.146775553749
$ws.Range("R9").Value = 4888.32097998374
$ws.Range("S9").Value = 0.1591775470895101
$ws.Range("T9").Value = 0.15917754708951

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 12.18925266666667
$ws.Range("H10").Value = 36.567758
$ws.Range("I10").Value = 0.04451691386950307
$ws.Range("J10").Value = 0.04451691386950307
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.397241999999999
$ws.Range("N10").Value = 16.191726
$ws.Range("O10").Value = 0.4331003391330405
$ws.Range("P10").Value = 0.4331003391330406
$ws.Range("Q10").Value = 65.78834644114532
$ws.Range("R10").Value = 592.0951179703079
$ws.Range("S10").Value = 0.01928029049403814
$ws.Range("T10").Value = 0.01928029049403814

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 12.18925266666667
$ws.Range("H11").Value = 36.567758
$ws.Range("I11").Value = 0.04451691386950307
$ws.Range("J11").Value = 0.04451691386950307
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.372979
$ws.Range("N11").Value = 13.118937
$ws.Range("O11").Value = 0.3509086099755513
$ws.Range("P11").Value = 0.3509086099755513
$ws.Range("Q11").Value = 53.30334593702733
$ws.Range("R11").Value = 479.7301134332459
$ws.Range("S11").Value = 0.01562136836634867
$ws.Range("T11").Value = 0.01562136836634866

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 12.18925266666667
$ws.Range("H12").Value = 36.567758
$ws.Range("I12").Value = 0.04451691386950307
$ws.Range("J12").Value = 0.04451691386950307
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.59901
$ws.Range("N12").Value = 1.79703
$ws.Range("O12").Value = 0.04806740815847847
$ws.Range("P12").Value = 0.04806740815847847
$ws.Range("Q12").Value = 7.301484239860001
$ws.Range("R12").Value = 65.71335815874001
$ws.Range("S12").Value = 0.002139812668921235
$ws.Range("T12").Value = 0.002139812668921235

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 12.18925266666667
$ws.Range("H13").Value = 36.567758
$ws.Range("I13").Value = 0.04451691386950307
$ws.Range("J13").Value = 0.04451691386950307
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.092643333333333
$ws.Range("N13").Value = 6.27793
$ws.Range("O13").Value = 0.1679236427329297
$ws.Range("P13").Value = 0.1679236427329297
$ws.Range("Q13").Value = 25.50775833121555
$ws.Range("R13").Value = 229.56982498094
$ws.Range("S13").Value = 0.007475442340195038
$ws.Range("T13").Value = 0.007475442340195038

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.5133286666666667
$ws.Range("H14").Value = 1.539986
$ws.Range("I14").Value = 0.001874750541781658
$ws.Range("J14").Value = 0.001874750541781658
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 5.397241999999999
$ws.Range("N14").Value = 16.191726
$ws.Range("O14").Value = 0.4331003391330405
$ws.Range("P14").Value = 0.4331003391330406
$ws.Range("Q14").Value = 2.770559039537333
$ws.Range("R14").Value = 24.935031355836
$ws.Range("S14").Value = 0.0008119550954354876
$ws.Range("T14").Value = 0.0008119550954354876

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.5133286666666667
$ws.Range("H15").Value = 1.539986
$ws.Range("I15").Value = 0.001874750541781658
$ws.Range("J15").Value = 0.001874750541781658
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.372979
$ws.Range("N15").Value = 13.118937
$ws.Range("O15").Value = 0.3509086099755513
$ws.Range("P15").Value = 0.3509086099755513
$ws.Range("Q15").Value = 2.244775479431333
$ws.Range("R15").Value = 20.202979314882
$ws.Range("S15").Value = 0.0006578661066675134
$ws.Range("T15").Value = 0.0006578661066675133

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.5133286666666667
$ws.Range("H16").Value = 1.539986
$ws.Range("I16").Value = 0.001874750541781658
$ws.Range("J16").Value = 0.001874750541781658
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.59901
$ws.Range("N16").Value = 1.79703
$ws.Range("O16").Value = 0.04806740815847847
$ws.Range("P16").Value = 0.04806740815847847
$ws.Range("Q16").Value = 0.30748900462
$ws.Range("R16").Value = 2.76740104158
$ws.Range("S16").Value = 0.0000901143994871476
$ws.Range("T16").Value = 0.0000901143994871476

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.5133286666666667
$ws.Range("H17").Value = 1.539986
$ws.Range("I17").Value = 0.001874750541781658
$ws.Range("J17").Value = 0.001874750541781658
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 2.092643333333333
$ws.Range("N17").Value = 6.27793
$ws.Range("O17").Value = 0.1679236427329297
$ws.Range("P17").Value = 0.1679236427329297
$ws.Range("Q17").Value = 1.074213812108889
$ws.Range("R17").Value = 9.667924308979998
$ws.Range("S17").Value = 0.0003148149401915096
$ws.Range("T17").Value = 0.0003148149401915096

Write-Host "done"